# Add a "Constraints" column (C) to the "Posts" sheet, populate it for the
# field-definition rows, make "Posts" the active/selected sheet (it was
# "comments" before), and leave the selection on the last populated cell.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Posts")

# New header cell, written first so the "Constraints" shared string is
# inserted before the constraint values themselves.
$ws1.Range("C1").Value = "Constraints"

# Copy the bold/filled header style from B1 onto the new header cell C1.
$null = $ws1.Range("B1").Copy()
$null = $ws1.Range("C1").PasteSpecial(-4122)

# Constraint values for each field row (Id, Title, Content, IsActive).
# Written in this order so the shared-string table fills up in the same
# sequence as the target workbook.
$ws1.Range("C2").Value = "['required', 'type' => 'number']"
$ws1.Range("C5").Value = "['required' => false]"
$ws1.Range("C3").Value = "['required' => true, 'maxlength' => '255']"
$ws1.Range("C4").Value = "['required' => true, 'type' => 'textarea', 'length' => ['min' => 10, 'max' => 512]]"

# CreatedAt (row 6) / UpdatedAt (row 7) intentionally have no constraints yet.

# Make "Posts" the active sheet/tab (previously "comments" was active) and
# leave the cursor on the newly-edited cell.
$null = $ws1.Activate()
$null = $ws1.Range("C5").Select()
